# "Nieuwe factuur en software update"
# Insert a new line item row (allekabels.nl invoice) above the existing
# "Bedrading, connectors, e.d." row, shifting everything below it down by
# one row (old row 19 -> 20, ..., old row 37 -> 38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a fresh row at position 19; Excel shifts formulas/styles/rows down
# and widens the shared-formula ranges (E2:E28 -> E2:E29, E30:E35 -> E31:E36)
# automatically, matching the diff.
$ws.Rows("19").Insert()

# Populate the new row 19 with the new invoice line item.
$ws.Range("B19").Value = "Bedrading, connectors, e.d."
$ws.Range("C19").Value = "allekabels.nl"
$ws.Range("F19").Value = "V"
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 22.93
$ws.Range("J19").Value = "Plastic-spray"

# Update the active selection to match the post-edit state.
$ws.Range("A19").Select()
